# Auto-generated Excel COM-interop script
# Applies targeted numeric cell updates (columns H-N) across sheets
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR as described by the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 353.61765
$ws.Range("J17").Value = 296.7742
$ws.Range("L17").Value = 890.3226
$ws.Range("N17").Value = -1226.3226
$ws.Range("H112").Value = 1295.1923
$ws.Range("J112").Value = 1323.9584
$ws.Range("L112").Value = 3971.8752
$ws.Range("N112").Value = -6187.8752
$ws.Range("H131").Value = 14349.214
$ws.Range("I131").Value = 12911.125
$ws.Range("K131").Value = 38733.375
$ws.Range("M131").Value = -33693.375
$ws.Range("H133").Value = 55750
$ws.Range("I133").Value = 30000
$ws.Range("J133").Value = 59428.57
$ws.Range("K133").Value = 30000
$ws.Range("L133").Value = 59428.57
$ws.Range("N133").Value = -69548.57000000001
$ws.Range("M133").Value = -24940
$ws.Range("H135").Value = 5274.909
$ws.Range("I135").Value = 4671
$ws.Range("J135").Value = 5999.6
$ws.Range("K135").Value = 42039
$ws.Range("L135").Value = 53996.4
$ws.Range("M135").Value = -39504
$ws.Range("N135").Value = -59066.4
$ws.Range("H138").Value = 2980.3262
$ws.Range("J138").Value = 3656.4546
$ws.Range("L138").Value = 10969.3638
$ws.Range("N138").Value = -21249.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2002.58
$ws.Range("I32").Value = 1867.75
$ws.Range("J32").Value = 2991.3333
$ws.Range("K32").Value = 1867.75
$ws.Range("L32").Value = 2991.3333
$ws.Range("N32").Value = -3565.3333
$ws.Range("M32").Value = -1580.75
$ws.Range("H74").Value = 4525.517
$ws.Range("I74").Value = 3930.7827
$ws.Range("J74").Value = 6805.3335
$ws.Range("K74").Value = 3930.7827
$ws.Range("L74").Value = 6805.3335
$ws.Range("M74").Value = -3056.7827
$ws.Range("N74").Value = -8553.333500000001
$ws.Range("H77").Value = 4525.517
$ws.Range("I77").Value = 3930.7827
$ws.Range("J77").Value = 6805.3335
$ws.Range("K77").Value = 19653.9135
$ws.Range("L77").Value = 34026.6675
$ws.Range("M77").Value = -15285.9135
$ws.Range("N77").Value = -42762.6675
$ws.Range("H132").Value = 3315
$ws.Range("I132").Value = 3315
$ws.Range("K132").Value = 9945
$ws.Range("M132").Value = -7415

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 33401778
$ws.Range("I86").Value = 50101668
$ws.Range("K86").Value = 50101668
$ws.Range("M86").Value = -50100545
$ws.Range("H89").Value = 33401778
$ws.Range("I89").Value = 50101668
$ws.Range("K89").Value = 250508340
$ws.Range("M89").Value = -250502724
$ws.Range("H107").Value = 5024.7095
$ws.Range("I107").Value = 5260.478
$ws.Range("K107").Value = 5260.478
$ws.Range("M107").Value = -3340.478
$ws.Range("H134").Value = 3047.7144
$ws.Range("I134").Value = 3047.7144
$ws.Range("K134").Value = 9143.143199999999
$ws.Range("M134").Value = -6608.143199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3021.3
$ws.Range("I31").Value = 1957.4615
$ws.Range("J31").Value = 4997
$ws.Range("K31").Value = 1957.4615
$ws.Range("L31").Value = 4997
$ws.Range("M31").Value = -1662.4615
$ws.Range("N31").Value = -5587
$ws.Range("H34").Value = 3021.3
$ws.Range("I34").Value = 1957.4615
$ws.Range("J34").Value = 4997
$ws.Range("K34").Value = 1957.4615
$ws.Range("L34").Value = 4997
$ws.Range("M34").Value = -1755.4615
$ws.Range("N34").Value = -5401
$ws.Range("H58").Value = 2620.2856
$ws.Range("I58").Value = 1451.2222
$ws.Range("K58").Value = 1451.2222
$ws.Range("M58").Value = -1248.2222
$ws.Range("H86").Value = 30651.941
$ws.Range("I86").Value = 36210.1
$ws.Range("J86").Value = 22711.715
$ws.Range("K86").Value = 36210.1
$ws.Range("L86").Value = 22711.715
$ws.Range("M86").Value = -35087.1
$ws.Range("N86").Value = -24957.715
$ws.Range("H89").Value = 30651.941
$ws.Range("I89").Value = 36210.1
$ws.Range("J89").Value = 22711.715
$ws.Range("K89").Value = 181050.5
$ws.Range("L89").Value = 113558.575
$ws.Range("M89").Value = -175434.5
$ws.Range("N89").Value = -124790.575
$ws.Range("H132").Value = 2982.1667
$ws.Range("J132").Value = 3742.8
$ws.Range("L132").Value = 11228.4
$ws.Range("N132").Value = -16288.4
$ws.Range("H134").Value = 5908.0713
$ws.Range("I134").Value = 6227.4165
$ws.Range("J134").Value = 3992
$ws.Range("K134").Value = 18682.2495
$ws.Range("L134").Value = 11976
$ws.Range("M134").Value = -16147.2495
$ws.Range("N134").Value = -17046
$ws.Range("H136").Value = 2620.2856
$ws.Range("I136").Value = 1451.2222
$ws.Range("K136").Value = 4353.6666
$ws.Range("M136").Value = -1803.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6999
$ws.Range("I56").Value = 6999
$ws.Range("K56").Value = 6999
$ws.Range("M56").Value = -6469
$ws.Range("H107").Value = 1957.6061
$ws.Range("I107").Value = 2344
$ws.Range("K107").Value = 7032
$ws.Range("M107").Value = -5112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 677666.7
$ws.Range("J21").Value = 8666.666999999999
$ws.Range("L21").Value = 8666.666999999999
$ws.Range("N21").Value = -9012.666999999999
$ws.Range("H30").Value = 677666.7
$ws.Range("J30").Value = 8666.666999999999
$ws.Range("L30").Value = 8666.666999999999
$ws.Range("N30").Value = -8876.666999999999
$ws.Range("H113").Value = 17317.727
$ws.Range("I113").Value = 4249.5
$ws.Range("J113").Value = 32999.6
$ws.Range("K113").Value = 4249.5
$ws.Range("L113").Value = 32999.6
$ws.Range("M113").Value = -2079.5
$ws.Range("N113").Value = -37339.6
$ws.Range("H120").Value = 49817.5
$ws.Range("J120").Value = 49817.5
$ws.Range("L120").Value = 49817.5
$ws.Range("N120").Value = -59493.5
$ws.Range("H122").Value = 6386.3335
$ws.Range("I122").Value = 6374.5
$ws.Range("J122").Value = 6399.857
$ws.Range("K122").Value = 19123.5
$ws.Range("L122").Value = 19199.571
$ws.Range("M122").Value = -16673.5
$ws.Range("N122").Value = -24099.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 20010
$ws.Range("J4").Value = 20010
$ws.Range("L4").Value = 20010
$ws.Range("N4").Value = -20236
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H21").Value = 1541.4
$ws.Range("I21").Value = 498.25
$ws.Range("J21").Value = 5714
$ws.Range("K21").Value = 498.25
$ws.Range("L21").Value = 5714
$ws.Range("M21").Value = -324.25
$ws.Range("N21").Value = -6062
$ws.Range("H22").Value = 11364873
$ws.Range("I22").Value = 22727898
$ws.Range("J22").Value = 1847.75
$ws.Range("K22").Value = 22727898
$ws.Range("L22").Value = 1847.75
$ws.Range("M22").Value = -22727603
$ws.Range("N22").Value = -2437.75
$ws.Range("H27").Value = 11364873
$ws.Range("I27").Value = 22727898
$ws.Range("J27").Value = 1847.75
$ws.Range("K27").Value = 22727898
$ws.Range("L27").Value = 1847.75
$ws.Range("M27").Value = -22727791
$ws.Range("N27").Value = -2061.75
$ws.Range("H28").Value = 20010
$ws.Range("J28").Value = 20010
$ws.Range("L28").Value = 20010
$ws.Range("N28").Value = -20474
$ws.Range("H37").Value = 20010
$ws.Range("J37").Value = 20010
$ws.Range("L37").Value = 20010
$ws.Range("N37").Value = -20224
$ws.Range("H55").Value = 499.90475
$ws.Range("I55").Value = 418.625
$ws.Range("J55").Value = 760
$ws.Range("K55").Value = 418.625
$ws.Range("L55").Value = 760
$ws.Range("M55").Value = -245.625
$ws.Range("N55").Value = -1106
$ws.Range("H61").Value = 2458.6365
$ws.Range("I61").Value = 2449.5557
$ws.Range("J61").Value = 2499.5
$ws.Range("K61").Value = 2449.5557
$ws.Range("L61").Value = 2499.5
$ws.Range("M61").Value = -2247.5557
$ws.Range("N61").Value = -2903.5
$ws.Range("H113").Value = 2458.6365
$ws.Range("I113").Value = 2449.5557
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 2449.5557
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = -279.5556999999999
$ws.Range("N113").Value = -6839.5
$ws.Range("H136").Value = 3534.125
$ws.Range("I136").Value = 3799.3
$ws.Range("J136").Value = 3092.1667
$ws.Range("K136").Value = 11397.9
$ws.Range("L136").Value = 9276.500100000001
$ws.Range("M136").Value = -8847.900000000001
$ws.Range("N136").Value = -14376.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 33744
$ws.Range("J95").Value = 33744
$ws.Range("L95").Value = 33744
$ws.Range("N95").Value = -39236
$ws.Range("H132").Value = 5240.8887
$ws.Range("I132").Value = 4029.7827
$ws.Range("K132").Value = 12089.3481
$ws.Range("M132").Value = -9559.348100000001
$ws.Range("H135").Value = 43374.125
$ws.Range("J135").Value = 43374.125
$ws.Range("L135").Value = 43374.125
$ws.Range("N135").Value = -53514.125
$ws.Range("H136").Value = 97603.8
$ws.Range("I136").Value = 107893.11
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 323679.33
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -321129.33
$ws.Range("N136").Value = -20100
